# Add a second greeting ("Waduppam!") as its own new paragraph right
# after the existing "Wassuppp!" paragraph.

$d = $word.ActiveDocument

# Find the existing "Wassuppp!" text and collapse the range to its end -
# i.e. the point right after the "!" and right before the (hidden)
# _GoBack bookmark that Word leaves at the end of that paragraph.
$r = $d.Content
$found = $r.Find.Execute("Wassuppp!", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0)
$r.Collapse(0)   # wdCollapseEnd

# Insert the new sentence's text right there (still inside the original
# paragraph, after "Wassuppp!" but before the trailing bookmark).
$r.InsertAfter("Waduppam!")

# Move back to the boundary between "Wassuppp!" and "Waduppam!" and split
# the paragraph there, so "Waduppam!" (together with the bookmark that
# was trailing it) becomes its own, new paragraph.
$r.Collapse(1)   # wdCollapseStart
$r.InsertParagraphAfter()
